# Update cryptos list values (Price/Volume) per upstream diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.649.42"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "1.637.63"
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.53"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.523"
$ws.Range("E6").Value = "  -1.55%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.97"
$ws.Range("E8").Value = "  -2.28%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0610"
$ws.Range("E10").Value = "  -0.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0893"
$ws.Range("E11").Value = "  +0.10%  "
$ws.Range("D12").Value = "1.870.17"
$ws.Range("E12").Value = "  -0.46%  "
$ws.Range("D13").Value = "1.627.59"
$ws.Range("E13").Value = "  -1.06%  "
$ws.Range("E15").Value = "  -5.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.60"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").Value = "27.640.62"
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.79"
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.77"
$ws.Range("E19").Value = "  +2.67%  "
$ws.Range("E20").Value = "  -0.38%  "
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("E22").Value = "  -0.95%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.07"
$ws.Range("E23").Value = "  +3.38%  "
$ws.Range("E24").Value = "  -1.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.50"
$ws.Range("E25").Value = "  +1.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.95"
$ws.Range("E26").Value = "  -1.26%  "
$ws.Range("E27").Value = "  -2.11%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("E29").Value = "  -0.29%  "
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0487"
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("E32").Value = "  -0.43%  "
$ws.Range("D33").Value = "1.458.92"
$ws.Range("E33").Value = "  +2.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.12"
$ws.Range("E34").Value = "  -2.06%  "
$ws.Range("E35").Value = "  -0.91%  "
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.563"
$ws.Range("E37").Value = "  -1.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.877"
$ws.Range("E38").Value = "  -1.44%  "
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.893"
$ws.Range("E40").Value = "  +9.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "70.12"
$ws.Range("E41").Value = "  +8.04%  "
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("E43").Value = "  -1.02%  "
$ws.Range("E44").Value = "  +1.15%  "
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("E46").Value = "  -0.53%  "
$ws.Range("D47").Value = "1.779.64"
$ws.Range("E47").Value = "  -0.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.73"
$ws.Range("E48").Value = "  +2.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.54"
$ws.Range("E49").Value = "  -2.13%  "
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0988"
$ws.Range("E51").Value = "  -0.45%  "
